$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a "pure" number (e.g. "313.84") need special handling:
# a plain .Value assignment would let Excel auto-convert the text into a
# numeric value (losing formatting such as trailing/leading zeros and the
# original inline-string type). Entering it as a formula that evaluates to a
# text string, then collapsing the formula to its static value via
# Copy/PasteSpecial(xlPasteValues), keeps the cell a plain text cell (matching
# the original t="inlineStr" cells) without disturbing any cell styles.
function Set-TextValue($addr, $text) {
    $cell = $ws.Range($addr)
    $escaped = $text.Replace("""", """""")
    $cell.Formula = "=""" + $escaped + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$ws.Range("D2").Value = "24.603.05"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.688.92"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.13%  "
Set-TextValue "D5" "313.84"
$ws.Range("E5").Value = "  -0.52%  "
Set-TextValue "D7" "0.3898"
$ws.Range("E7").Value = "  -1.15%  "
$ws.Range("E8").Value = "  -0.56%  "
Set-TextValue "D9" "1.497"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("E10").Value = "  +0.16%  "
Set-TextValue "D11" "52.65"
$ws.Range("E11").Value = "  +0.35%  "
Set-TextValue "D12" "0.08753"
$ws.Range("E12").Value = "  -0.95%  "
Set-TextValue "D13" "7.569"
$ws.Range("E13").Value = "  +4.48%  "
Set-TextValue "D14" "24.90"
$ws.Range("E14").Value = "  +5.94%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D15" "0.00001350"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D16" "7.958"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "1.681.68"
$ws.Range("E17").Value = "  -0.58%  "
Set-TextValue "D18" "98.52"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("E19").Value = "  +1.31%  "
Set-TextValue "D20" "19.86"
$ws.Range("E20").Value = "  +1.79%  "
Set-TextValue "D21" "7.284"
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("E22").Value = "  +0.01%  "
Set-TextValue "D23" "14.25"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "24.587.94"
$ws.Range("E24").Value = "  -0.26%  "
Set-TextValue "D25" "3.005"
$ws.Range("E25").Value = "  -9.16%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.14%  "
Set-TextValue "D28" "161.92"
$ws.Range("E28").Value = "  -0.39%  "
Set-TextValue "D29" "8.748"
$ws.Range("E29").Value = "  +15.17%  "
Set-TextValue "D30" "136.76"
$ws.Range("E30").Value = "  +0.97%  "
Set-TextValue "D31" "5.218"
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").Value = "1.866.40"
$ws.Range("E32").Value = "  -0.67%  "
Set-TextValue "D33" "0.08864"
$ws.Range("E33").Value = "  +3.77%  "
Set-TextValue "D34" "7.406"
$ws.Range("E34").Value = "  +4.46%  "
Set-TextValue "D35" "1.039"
$ws.Range("E35").Value = "  -1.68%  "
Set-TextValue "D36" "1.984"
$ws.Range("E36").Value = "  +5.30%  "
Set-TextValue "D37" "0.02924"
$ws.Range("E37").Value = "  +7.52%  "
Set-TextValue "D38" "0.2734"
$ws.Range("E38").Value = "  -0.06%  "
Set-TextValue "D39" "10.77"
Set-TextValue "D40" "0.09139"
$ws.Range("E40").Value = "  -0.59%  "
Set-TextValue "D41" "14.18"
$ws.Range("E41").Value = "  -1.87%  "
Set-TextValue "D42" "0.7816"
$ws.Range("E42").Value = "  +2.45%  "
Set-TextValue "D43" "1.459"
$ws.Range("E43").Value = "  -0.27%  "
Set-TextValue "D44" "16.56"
$ws.Range("E44").Value = "  +2.67%  "
Set-TextValue "D45" "0.7189"
$ws.Range("E45").Value = "  +0.82%  "
Set-TextValue "D46" "2.592"
$ws.Range("E46").Value = "  +0.18%  "
Set-TextValue "D47" "4.193"
$ws.Range("E47").Value = "  -0.31%  "
Set-TextValue "D49" "1.337"
$ws.Range("E49").Value = "  +1.47%  "
Set-TextValue "D50" "137.78"
$ws.Range("E50").Value = "  -1.55%  "
Set-TextValue "D51" "90.98"
$ws.Range("E51").Value = "  +1.29%  "
